# PLTR model update:
#  - bump the near-term revenue growth multiplier row (Sheet2 row 2, cols O:W)
#    from a flat 34% YoY ramp to 50% in the first forecast year then 40%
#    thereafter
#  - trim the long-run discount rate (Z20) from 10% to 9.5% and show it (and
#    the two rates above it) with fewer/standardised decimal places
#  - leave the active selection on Z18 to match where the author ended up

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Row 2: revenue build multipliers --------------------------------------
$ws.Range("O2").Formula = "=N2*1.5"
$ws.Range("P2").Formula = "=O2*1.4"
$ws.Range("Q2").Formula = "=P2*1.4"
$ws.Range("R2").Formula = "=Q2*1.4"
$ws.Range("S2").Formula = "=R2*1.4"
$ws.Range("T2:W2").Formula = "=S2*1.4"

# --- Assumptions block (Y18:Z20) --------------------------------------------
# Growth-rate inputs: same values, displayed as whole-percent instead of
# two-decimal percent.
$ws.Range("Z18").NumberFormat = "0%"
$ws.Range("Z19").NumberFormat = "0%"

# Discount rate: 10.0% -> 9.5%, shown to one decimal place.
$ws.Range("Z20").NumberFormat = "0.0%"
$ws.Range("Z20").Value = 0.095

# --- Leave the user's selection where the edit ended up ---------------------
[void]$ws.Range("Z18").Select()
